$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.855.52'
$ws.Range("E2").Value = '  +3.63%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.708.21'
$ws.Range("E3").Value = '  +3.17%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.11%  '

# Row 5
$ws.Range("E5").Value = '  +1.41%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.26'
$ws.Range("E6").Value = '  +1.41%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.11%  '

# Row 8
$ws.Range("E8").Value = '  +1.62%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.721.65'
$ws.Range("E9").Value = '  +3.44%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.05'
$ws.Range("E10").Value = '  +12.07%  '

# Row 11
$ws.Range("E11").Value = '  +1.69%  '

# Row 12
$ws.Range("E12").Value = '  +2.13%  '

# Row 13
$ws.Range("E13").Value = '  +3.42%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.186.52'
$ws.Range("E14").Value = '  +3.30%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.801.26'
$ws.Range("E15").Value = '  +3.52%  '

# Row 16
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.37'
$ws.Range("E16").Value = '  +3.14%  '

# Row 17
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.726.27'
$ws.Range("E17").Value = '  +3.64%  '

# Row 18
$ws.Range("E18").Value = '  +1.63%  '

# Row 19
$ws.Range("E19").Value = '  +2.20%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '344.06'
$ws.Range("E20").Value = '  -0.23%  '

# Row 21
$ws.Range("E21").Value = '  +3.86%  '

# Row 22
$ws.Range("E22").Value = '  +5.47%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.12%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.73'
$ws.Range("E24").Value = '  +3.81%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.419'
$ws.Range("E25").Value = '  +1.36%  '

# Row 26
$ws.Range("E26").Value = '  +4.07%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.996'
$ws.Range("E27").Value = '  -0.06%  '

# Row 29
$ws.Range("E29").Value = '  +3.49%  '

# Row 30
$ws.Range("E30").Value = '  +9.59%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  +0.04%  '

# Row 32
$ws.Range("E32").Value = '  +2.12%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.01'
$ws.Range("E33").Value = '  +1.05%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '149.67'
$ws.Range("E34").Value = '  -0.22%  '

# Row 35
$ws.Range("E35").Value = '  +7.96%  '

# Row 36
$ws.Range("E36").Value = '  +8.64%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.928'
$ws.Range("E37").Value = '  -4.85%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.879'
$ws.Range("E38").Value = '  +5.43%  '

# Row 39
$ws.Range("E39").Value = '  +7.57%  '

# Row 40
$ws.Range("E40").Value = '  +1.57%  '

# Row 41
$ws.Range("E41").Value = '  +1.38%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '281.07'
$ws.Range("E42").Value = '  +1.52%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.14'
$ws.Range("E43").Value = '  +3.50%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  +0.16%  '

# Row 45
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0988'
$ws.Range("E45").Value = '  +0.75%  '

# Row 46
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.611'
$ws.Range("E46").Value = '  +2.20%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.129.55'
$ws.Range("E47").Value = '  +7.13%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0543'
$ws.Range("E48").Value = '  +4.61%  '

# Row 49
$ws.Range("E49").Value = '  +5.87%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '10.54'
$ws.Range("E50").Value = '  +2.25%  '

# Row 51
$ws.Range("E51").Value = '  +1.94%  '

